$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("D2").Value = '59.993.52'
$ws.Range("E2").Value = '  -6.45%  '
$ws.Range("D3").Value = '3.294.66'
$ws.Range("E3").Value = '  -5.28%  '
$ws.Range("E4").Value = '  -0.04%  '
$ws.Range("D5").Value = "'562.68"
$ws.Range("E5").Value = '  -3.90%  '
$ws.Range("D6").Value = "'128.32"
$ws.Range("E6").Value = '  -2.83%  '
$ws.Range("E7").Value = '  -0.05%  '
$ws.Range("D8").Value = '3.296.85'
$ws.Range("E8").Value = '  -5.16%  '
$ws.Range("E9").Value = '  -2.55%  '
$ws.Range("D10").Value = "'7.42"
$ws.Range("E10").Value = '  -3.89%  '
$ws.Range("E11").Value = '  -5.82%  '
$ws.Range("E12").Value = '  -4.11%  '
$ws.Range("D13").Value = '3.855.78'
$ws.Range("E13").Value = '  -5.32%  '
$ws.Range("E14").Value = '  -0.40%  '
$ws.Range("D15").Value = '3.296.13'
$ws.Range("E15").Value = '  -5.28%  '
$ws.Range("E16").Value = '  -6.58%  '
$ws.Range("D17").Value = '60.202.77'
$ws.Range("E17").Value = '  -6.14%  '
$ws.Range("D18").Value = "'24.00"
$ws.Range("E18").Value = '  -4.51%  '
$ws.Range("D19").Value = "'5.59"
$ws.Range("E19").Value = '  -1.78%  '
$ws.Range("D20").Value = "'13.20"
$ws.Range("E20").Value = '  -1.46%  '
$ws.Range("D21").Value = "'8.85"
$ws.Range("D22").Value = "'348.08"
$ws.Range("E22").Value = '  -9.83%  '
$ws.Range("D23").Value = "'0.551"
$ws.Range("E23").Value = '  -2.75%  '
$ws.Range("E24").Value = '  +0.06%  '
$ws.Range("D25").Value = '3.423.91'
$ws.Range("E25").Value = '  -5.42%  '
$ws.Range("D26").Value = "'68.65"
$ws.Range("E26").Value = '  -7.94%  '
$ws.Range("E27").Value = '  -4.30%  '
$ws.Range("D28").Value = "'0.994"
$ws.Range("E28").Value = '  -0.56%  '
$ws.Range("D29").Value = "'7.26"
$ws.Range("E29").Value = '  +2.31%  '
$ws.Range("E30").Value = '  +0.72%  '
$ws.Range("E31").Value = '  -2.63%  '
$ws.Range("E32").Value = '  -2.27%  '
$ws.Range("E33").Value = '  -6.02%  '
$ws.Range("E34").Value = '  +0.00%  '
$ws.Range("D35").Value = '3.323.38'
$ws.Range("D36").Value = "'22.55"
$ws.Range("E36").Value = '  -1.83%  '
$ws.Range("D37").Value = "'5.30"
$ws.Range("E37").Value = '  +1.39%  '
$ws.Range("E38").Value = '  -0.67%  '
$ws.Range("E39").Value = '  -1.90%  '
$ws.Range("D40").Value = "'156.46"
$ws.Range("E40").Value = '  -3.63%  '
$ws.Range("E41").Value = '  -4.11%  '
$ws.Range("E42").Value = '  -0.04%  '
$ws.Range("D43").Value = "'40.62"
$ws.Range("E43").Value = '  -2.20%  '
$ws.Range("D44").Value = "'4.29"
$ws.Range("E44").Value = '  -1.12%  '
$ws.Range("E45").Value = '  -7.39%  '
$ws.Range("E46").Value = '  +2.71%  '
$ws.Range("D47").Value = "'22.47"
$ws.Range("E47").Value = '  -4.46%  '
$ws.Range("E48").Value = '  -5.66%  '
$ws.Range("D49").Value = "'6.68"
$ws.Range("E49").Value = '  -0.61%  '
$ws.Range("E50").Value = '  +6.08%  '
$ws.Range("D51").Value = "'0.848"
$ws.Range("E51").Value = '  -5.91%  '
